$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rearrange match data within same-date blocks (columns F:V) ---
# Values captured from the source workbook (see commit diff) and written
# directly to their new row positions; column A (index) and E (date) stay put.

# row 36  (was row 37: Chelmianka Chelm vs Podhale Nowy Targ)
$ws.Cells.Item(36,6).Value = 'Chelmianka Chelm'
$ws.Cells.Item(36,7).Value = 2
$ws.Cells.Item(36,8).Value = 'Podhale Nowy Targ'
$ws.Cells.Item(36,9).Value = 0
$ws.Cells.Item(36,10).Value = 2.47
$ws.Cells.Item(36,11).Value = '26/08/2023 05:12'
$ws.Cells.Item(36,12).Value = 2.47
$ws.Cells.Item(36,13).Value = '26/08/2023 05:12'
$ws.Cells.Item(36,14).Value = 3.14
$ws.Cells.Item(36,15).Value = '26/08/2023 05:12'
$ws.Cells.Item(36,16).Value = 3.24
$ws.Cells.Item(36,17).Value = '27/08/2023 15:04'
$ws.Cells.Item(36,18).Value = 2.36
$ws.Cells.Item(36,19).Value = '26/08/2023 05:12'
$ws.Cells.Item(36,20).Value = 2.36
$ws.Cells.Item(36,21).Value = '26/08/2023 05:12'
$ws.Cells.Item(36,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/chelmianka-chelm-podhale-nowy-targ/rVthZBJo/'

# row 37  (was row 38: Avia Swidnik vs KS Wieczysta Krakow)
$ws.Cells.Item(37,6).Value = 'Avia Swidnik'
$ws.Cells.Item(37,7).Value = 3
$ws.Cells.Item(37,8).Value = 'KS Wieczysta Krakow'
$ws.Cells.Item(37,9).Value = 1
$ws.Cells.Item(37,10).Value = 3.03
$ws.Cells.Item(37,11).Value = '26/08/2023 05:12'
$ws.Cells.Item(37,12).Value = 3.05
$ws.Cells.Item(37,13).Value = '27/08/2023 16:58'
$ws.Cells.Item(37,14).Value = 3.29
$ws.Cells.Item(37,15).Value = '26/08/2023 05:12'
$ws.Cells.Item(37,16).Value = 3.7
$ws.Cells.Item(37,17).Value = '27/08/2023 16:36'
$ws.Cells.Item(37,18).Value = 1.99
$ws.Cells.Item(37,19).Value = '26/08/2023 05:12'
$ws.Cells.Item(37,20).Value = 1.8
$ws.Cells.Item(37,21).Value = '27/08/2023 16:58'
$ws.Cells.Item(37,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/avia-swidnik-ks-wieczysta-krakow/dYpdYVYi/'

# row 38  (was row 36: Orleta Radzyn vs Ostrowiec Swietokrzyski)
$ws.Cells.Item(38,6).Value = 'Orleta Radzyn'
$ws.Cells.Item(38,7).Value = 0
$ws.Cells.Item(38,8).Value = 'Ostrowiec Swietokrzyski'
$ws.Cells.Item(38,9).Value = 2
$ws.Cells.Item(38,10).Value = 2.96
$ws.Cells.Item(38,11).Value = '26/08/2023 05:12'
$ws.Cells.Item(38,12).Value = 2.7
$ws.Cells.Item(38,13).Value = '27/08/2023 16:11'
$ws.Cells.Item(38,14).Value = 3.31
$ws.Cells.Item(38,15).Value = '26/08/2023 05:12'
$ws.Cells.Item(38,16).Value = 3.38
$ws.Cells.Item(38,17).Value = '27/08/2023 16:11'
$ws.Cells.Item(38,18).Value = 1.97
$ws.Cells.Item(38,19).Value = '26/08/2023 05:12'
$ws.Cells.Item(38,20).Value = 2.28
$ws.Cells.Item(38,21).Value = '27/08/2023 16:11'
$ws.Cells.Item(38,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/orleta-radzyn-ostrowiec-swietokrzyski/Ac01VTIA/'

# row 40  (was row 41: KS Wieczysta Krakow vs Wiazownica)
$ws.Cells.Item(40,6).Value = 'KS Wieczysta Krakow'
$ws.Cells.Item(40,7).Value = 7
$ws.Cells.Item(40,8).Value = 'Wiazownica'
$ws.Cells.Item(40,9).Value = 0
$ws.Cells.Item(40,10).Value = 1.17
$ws.Cells.Item(40,11).Value = '01/09/2023 00:12'
$ws.Cells.Item(40,12).Value = 1.16
$ws.Cells.Item(40,13).Value = '02/09/2023 11:57'
$ws.Cells.Item(40,14).Value = 6.18
$ws.Cells.Item(40,15).Value = '01/09/2023 00:12'
$ws.Cells.Item(40,16).Value = 5.82
$ws.Cells.Item(40,17).Value = '02/09/2023 11:57'
$ws.Cells.Item(40,18).Value = 7.89
$ws.Cells.Item(40,19).Value = '01/09/2023 00:12'
$ws.Cells.Item(40,20).Value = 16.35
$ws.Cells.Item(40,21).Value = '02/09/2023 11:57'
$ws.Cells.Item(40,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/ks-wieczysta-krakow-wiazownica/2eutOoBp/'

# row 41  (was row 40: Czarni Polaniec vs Wislanie Jaskowice)
$ws.Cells.Item(41,6).Value = 'Czarni Polaniec'
$ws.Cells.Item(41,7).Value = 2
$ws.Cells.Item(41,8).Value = 'Wislanie Jaskowice'
$ws.Cells.Item(41,9).Value = 1
$ws.Cells.Item(41,10).Value = 1.9
$ws.Cells.Item(41,11).Value = '02/09/2023 04:12'
$ws.Cells.Item(41,12).Value = 1.92
$ws.Cells.Item(41,13).Value = '02/09/2023 11:57'
$ws.Cells.Item(41,14).Value = 3.55
$ws.Cells.Item(41,15).Value = '02/09/2023 04:12'
$ws.Cells.Item(41,16).Value = 3.94
$ws.Cells.Item(41,17).Value = '02/09/2023 11:57'
$ws.Cells.Item(41,18).Value = 3.18
$ws.Cells.Item(41,19).Value = '02/09/2023 04:12'
$ws.Cells.Item(41,20).Value = 3.04
$ws.Cells.Item(41,21).Value = '02/09/2023 11:57'
$ws.Cells.Item(41,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/czarni-polaniec-wislanie-jaskowice/jHLzpTX3/'

# row 42  (was row 45: Unia Tarnow vs Orleta Radzyn)
$ws.Cells.Item(42,6).Value = 'Unia Tarnow'
$ws.Cells.Item(42,7).Value = 2
$ws.Cells.Item(42,8).Value = 'Orleta Radzyn'
$ws.Cells.Item(42,9).Value = 0
$ws.Cells.Item(42,10).Value = 1.83
$ws.Cells.Item(42,11).Value = '01/09/2023 05:13'
$ws.Cells.Item(42,12).Value = 1.72
$ws.Cells.Item(42,13).Value = '02/09/2023 16:59'
$ws.Cells.Item(42,14).Value = 3.42
$ws.Cells.Item(42,15).Value = '01/09/2023 05:13'
$ws.Cells.Item(42,16).Value = 3.74
$ws.Cells.Item(42,17).Value = '02/09/2023 16:59'
$ws.Cells.Item(42,18).Value = 3.21
$ws.Cells.Item(42,19).Value = '01/09/2023 05:13'
$ws.Cells.Item(42,20).Value = 3.88
$ws.Cells.Item(42,21).Value = '02/09/2023 16:59'
$ws.Cells.Item(42,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/unia-tarnow-orleta-radzyn/rBb9T7mN/'

# row 43  (was row 44: Sokol Sieniawa vs Ostrowiec Swietokrzyski)
$ws.Cells.Item(43,6).Value = 'Sokol Sieniawa'
$ws.Cells.Item(43,7).Value = 1
$ws.Cells.Item(43,8).Value = 'Ostrowiec Swietokrzyski'
$ws.Cells.Item(43,9).Value = 3
$ws.Cells.Item(43,10).Value = 2.67
$ws.Cells.Item(43,11).Value = '01/09/2023 05:13'
$ws.Cells.Item(43,12).Value = 3.91
$ws.Cells.Item(43,13).Value = '02/09/2023 16:51'
$ws.Cells.Item(43,14).Value = 3.17
$ws.Cells.Item(43,15).Value = '01/09/2023 05:13'
$ws.Cells.Item(43,16).Value = 3.5
$ws.Cells.Item(43,17).Value = '02/09/2023 16:51'
$ws.Cells.Item(43,18).Value = 2.19
$ws.Cells.Item(43,19).Value = '01/09/2023 05:13'
$ws.Cells.Item(43,20).Value = 1.78
$ws.Cells.Item(43,21).Value = '02/09/2023 16:51'
$ws.Cells.Item(43,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/sokol-sieniawa-ostrowiec-swietokrzyski/lla5UmYG/'

# row 44  (was row 43: Siarka Tarnobrzeg vs Garbarnia)
$ws.Cells.Item(44,6).Value = 'Siarka Tarnobrzeg'
$ws.Cells.Item(44,7).Value = 0
$ws.Cells.Item(44,8).Value = 'Garbarnia'
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 2.23
$ws.Cells.Item(44,11).Value = '02/09/2023 13:42'
$ws.Cells.Item(44,12).Value = 1.98
$ws.Cells.Item(44,13).Value = '02/09/2023 16:32'
$ws.Cells.Item(44,14).Value = 3.35
$ws.Cells.Item(44,15).Value = '02/09/2023 13:42'
$ws.Cells.Item(44,16).Value = 3.3
$ws.Cells.Item(44,17).Value = '02/09/2023 16:35'
$ws.Cells.Item(44,18).Value = 2.68
$ws.Cells.Item(44,19).Value = '02/09/2023 13:42'
$ws.Cells.Item(44,20).Value = 2.99
$ws.Cells.Item(44,21).Value = '02/09/2023 16:32'
$ws.Cells.Item(44,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/siarka-tarnobrzeg-garbarnia/8hWVp9Ic/'

# row 45  (was row 42: Podhale Nowy Targ vs Avia Swidnik)
$ws.Cells.Item(45,6).Value = 'Podhale Nowy Targ'
$ws.Cells.Item(45,7).Value = 1
$ws.Cells.Item(45,8).Value = 'Avia Swidnik'
$ws.Cells.Item(45,9).Value = 6
$ws.Cells.Item(45,10).Value = 3
$ws.Cells.Item(45,11).Value = '01/09/2023 05:13'
$ws.Cells.Item(45,12).Value = 2.85
$ws.Cells.Item(45,13).Value = '02/09/2023 16:37'
$ws.Cells.Item(45,14).Value = 3.26
$ws.Cells.Item(45,15).Value = '01/09/2023 05:13'
$ws.Cells.Item(45,16).Value = 3.61
$ws.Cells.Item(45,17).Value = '02/09/2023 16:37'
$ws.Cells.Item(45,18).Value = 1.97
$ws.Cells.Item(45,19).Value = '01/09/2023 05:13'
$ws.Cells.Item(45,20).Value = 2.04
$ws.Cells.Item(45,21).Value = '02/09/2023 16:37'
$ws.Cells.Item(45,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/podhale-nowy-targ-avia-swidnik/tnvpN5Qj/'

# row 48  (was row 50: Avia Swidnik vs Swidniczanka Swidnik)
$ws.Cells.Item(48,6).Value = 'Avia Swidnik'
$ws.Cells.Item(48,7).Value = 1
$ws.Cells.Item(48,8).Value = 'Swidniczanka Swidnik'
$ws.Cells.Item(48,9).Value = 3
$ws.Cells.Item(48,10).Value = 1.12
$ws.Cells.Item(48,11).Value = '09/09/2023 12:43'
$ws.Cells.Item(48,12).Value = 1.27
$ws.Cells.Item(48,13).Value = '09/09/2023 15:31'
$ws.Cells.Item(48,14).Value = 7.56
$ws.Cells.Item(48,15).Value = '09/09/2023 12:43'
$ws.Cells.Item(48,16).Value = 7.25
$ws.Cells.Item(48,17).Value = '09/09/2023 15:31'
$ws.Cells.Item(48,18).Value = 11.01
$ws.Cells.Item(48,19).Value = '09/09/2023 12:43'
$ws.Cells.Item(48,20).Value = 5.59
$ws.Cells.Item(48,21).Value = '09/09/2023 15:31'
$ws.Cells.Item(48,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/avia-swidnik-swidniczanka-swidnik/lC3NzPfj/'

# row 49  (was row 51: Garbarnia vs Czarni Polaniec)
$ws.Cells.Item(49,6).Value = 'Garbarnia'
$ws.Cells.Item(49,7).Value = 0
$ws.Cells.Item(49,8).Value = 'Czarni Polaniec'
$ws.Cells.Item(49,9).Value = 0
$ws.Cells.Item(49,10).Value = 1.57
$ws.Cells.Item(49,11).Value = '09/09/2023 12:43'
$ws.Cells.Item(49,12).Value = 1.74
$ws.Cells.Item(49,13).Value = '09/09/2023 15:43'
$ws.Cells.Item(49,14).Value = 3.93
$ws.Cells.Item(49,15).Value = '09/09/2023 12:43'
$ws.Cells.Item(49,16).Value = 3.9
$ws.Cells.Item(49,17).Value = '09/09/2023 15:43'
$ws.Cells.Item(49,18).Value = 4.28
$ws.Cells.Item(49,19).Value = '09/09/2023 12:43'
$ws.Cells.Item(49,20).Value = 3.66
$ws.Cells.Item(49,21).Value = '09/09/2023 15:43'
$ws.Cells.Item(49,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/garbarnia-czarni-polaniec/htIitoXS/'

# row 50  (was row 49: Orleta Radzyn vs Wisloka Debica)
$ws.Cells.Item(50,6).Value = 'Orleta Radzyn'
$ws.Cells.Item(50,7).Value = 1
$ws.Cells.Item(50,8).Value = 'Wisloka Debica'
$ws.Cells.Item(50,9).Value = 1
$ws.Cells.Item(50,10).Value = 2.64
$ws.Cells.Item(50,11).Value = '08/09/2023 04:13'
$ws.Cells.Item(50,12).Value = 2.63
$ws.Cells.Item(50,13).Value = '09/09/2023 15:53'
$ws.Cells.Item(50,14).Value = 3.28
$ws.Cells.Item(50,15).Value = '08/09/2023 04:13'
$ws.Cells.Item(50,16).Value = 3.46
$ws.Cells.Item(50,17).Value = '09/09/2023 15:53'
$ws.Cells.Item(50,18).Value = 2.16
$ws.Cells.Item(50,19).Value = '08/09/2023 04:13'
$ws.Cells.Item(50,20).Value = 2.3
$ws.Cells.Item(50,21).Value = '09/09/2023 15:53'
$ws.Cells.Item(50,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/orleta-radzyn-wisloka-debica/rBpEcn23/'

# row 51  (was row 48: Wislanie Jaskowice vs Star Starachowice)
$ws.Cells.Item(51,6).Value = 'Wislanie Jaskowice'
$ws.Cells.Item(51,7).Value = 1
$ws.Cells.Item(51,8).Value = 'Star Starachowice'
$ws.Cells.Item(51,9).Value = 1
$ws.Cells.Item(51,10).Value = 2.57
$ws.Cells.Item(51,11).Value = '09/09/2023 12:43'
$ws.Cells.Item(51,12).Value = 2.68
$ws.Cells.Item(51,13).Value = '09/09/2023 15:53'
$ws.Cells.Item(51,14).Value = 3.33
$ws.Cells.Item(51,15).Value = '09/09/2023 12:43'
$ws.Cells.Item(51,16).Value = 3.42
$ws.Cells.Item(51,17).Value = '09/09/2023 15:53'
$ws.Cells.Item(51,18).Value = 2.32
$ws.Cells.Item(51,19).Value = '09/09/2023 12:43'
$ws.Cells.Item(51,20).Value = 2.28
$ws.Cells.Item(51,21).Value = '09/09/2023 15:53'
$ws.Cells.Item(51,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/wislanie-jaskowice-star-starachowice/z5NnsRHM/'

# --- Step 2: append two new rows (63, 64) ---
# Copy formatting from the last existing row (62) for the A (index) and E (date) styles

$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("E62").Copy()
$ws.Range("E63").PasteSpecial(-4122)
$ws.Cells.Item(63,1).Value = 62
$ws.Cells.Item(63,2).Value = 'poland'
$ws.Cells.Item(63,3).Value = 'iii-liga-group-iv'
$ws.Cells.Item(63,4).Value = '2023-2024'
$ws.Cells.Item(63,5).Value = 45191.66666666666
$ws.Cells.Item(63,6).Value = 'Garbarnia'
$ws.Cells.Item(63,7).Value = 2
$ws.Cells.Item(63,8).Value = 'Karpaty Krosno'
$ws.Cells.Item(63,9).Value = 1
$ws.Cells.Item(63,10).Value = 1.35
$ws.Cells.Item(63,11).Value = '22/09/2023 10:12'
$ws.Cells.Item(63,12).Value = 1.31
$ws.Cells.Item(63,13).Value = '22/09/2023 15:44'
$ws.Cells.Item(63,14).Value = 4.58
$ws.Cells.Item(63,15).Value = '22/09/2023 10:12'
$ws.Cells.Item(63,16).Value = 4.96
$ws.Cells.Item(63,17).Value = '22/09/2023 15:44'
$ws.Cells.Item(63,18).Value = 6.07
$ws.Cells.Item(63,19).Value = '22/09/2023 10:12'
$ws.Cells.Item(63,20).Value = 7.14
$ws.Cells.Item(63,21).Value = '22/09/2023 15:44'
$ws.Cells.Item(63,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/garbarnia-ks-karpaty-krosno/jeC383GL/'

$ws.Range("A62").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("E62").Copy()
$ws.Range("E64").PasteSpecial(-4122)
$ws.Cells.Item(64,1).Value = 63
$ws.Cells.Item(64,2).Value = 'poland'
$ws.Cells.Item(64,3).Value = 'iii-liga-group-iv'
$ws.Cells.Item(64,4).Value = '2023-2024'
$ws.Cells.Item(64,5).Value = 45191.81180555555
$ws.Cells.Item(64,6).Value = 'Ostrowiec Swietokrzyski'
$ws.Cells.Item(64,7).Value = 0
$ws.Cells.Item(64,8).Value = 'KS Wieczysta Krakow'
$ws.Cells.Item(64,9).Value = 1
$ws.Cells.Item(64,10).Value = 4.69
$ws.Cells.Item(64,11).Value = '21/09/2023 06:42'
$ws.Cells.Item(64,12).Value = 3.15
$ws.Cells.Item(64,13).Value = '22/09/2023 19:26'
$ws.Cells.Item(64,14).Value = 4.12
$ws.Cells.Item(64,15).Value = '21/09/2023 06:42'
$ws.Cells.Item(64,16).Value = 4.74
$ws.Cells.Item(64,17).Value = '22/09/2023 19:26'
$ws.Cells.Item(64,18).Value = 1.44
$ws.Cells.Item(64,19).Value = '21/09/2023 06:42'
$ws.Cells.Item(64,20).Value = 1.74
$ws.Cells.Item(64,21).Value = '22/09/2023 19:26'
$ws.Cells.Item(64,22).Value = 'https://www.betexplorer.com/football/poland/iii-liga-group-iv/ostrowiec-swietokrzyski-ks-wieczysta-krakow/jJfjaaw8/'

$excel.CutCopyMode = 0

"done"